$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1666.6666
$ws.Range("J2").Value = 2000
$ws.Range("L2").Value = 2000
$ws.Range("N2").Value = -2226
$ws.Range("H13").Value = 2500
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 2500
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 2500
$ws.Range("M13").Value = $null
$ws.Range("N13").Value = -2838
$ws.Range("H17").Value = 507.14285
$ws.Range("J17").Value = 550
$ws.Range("L17").Value = 1650
$ws.Range("N17").Value = -1986
$ws.Range("H29").Value = 2197.5
$ws.Range("I29").Value = 516
$ws.Range("J29").Value = 5000
$ws.Range("K29").Value = 1548
$ws.Range("L29").Value = 15000
$ws.Range("M29").Value = -1267
$ws.Range("N29").Value = -15562
$ws.Range("H58").Value = 574.75
$ws.Range("I58").Value = 574.75
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1724.25
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -1574.25
$ws.Range("N58").Value = $null
$ws.Range("H62").Value = 1458.75
$ws.Range("I62").Value = 1444.5
$ws.Range("K62").Value = 1444.5
$ws.Range("M62").Value = -820.5
$ws.Range("H65").Value = 1458.75
$ws.Range("I65").Value = 1444.5
$ws.Range("K65").Value = 7222.5
$ws.Range("M65").Value = -4102.5
$ws.Range("H132").Value = 3106.5
$ws.Range("I132").Value = 3193
$ws.Range("J132").Value = 2976.75
$ws.Range("K132").Value = 9579
$ws.Range("L132").Value = 8930.25
$ws.Range("M132").Value = -7049
$ws.Range("N132").Value = -13990.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 6511.3335
$ws.Range("I37").Value = 3022.6667
$ws.Range("J37").Value = 10000
$ws.Range("K37").Value = 3022.6667
$ws.Range("L37").Value = 10000
$ws.Range("M37").Value = -2749.6667
$ws.Range("N37").Value = -10546
$ws.Range("H44").Value = 20000
$ws.Range("J44").Value = 20000
$ws.Range("L44").Value = 20000
$ws.Range("N44").Value = -20976
$ws.Range("H63").Value = 2738.75
$ws.Range("J63").Value = 2178
$ws.Range("L63").Value = 2178
$ws.Range("N63").Value = -3550
$ws.Range("H66").Value = 2738.75
$ws.Range("J66").Value = 2178
$ws.Range("L66").Value = 10890
$ws.Range("N66").Value = -17754
$ws.Range("H74").Value = 3450
$ws.Range("I74").Value = 1900
$ws.Range("K74").Value = 1900
$ws.Range("M74").Value = -1026
$ws.Range("H77").Value = 3450
$ws.Range("I77").Value = 1900
$ws.Range("K77").Value = 9500
$ws.Range("M77").Value = -5132

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").Value = $null
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").Value = $null
$ws.Range("H134").Value = 1569.125
$ws.Range("I134").Value = 1522.5
$ws.Range("K134").Value = 4567.5
$ws.Range("M134").Value = -2032.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 20016.25
$ws.Range("J41").Value = 20016.25
$ws.Range("L41").Value = 20016.25
$ws.Range("N41").Value = -20872.25
$ws.Range("H50").Value = 23571.428
$ws.Range("J50").Value = 23571.428
$ws.Range("L50").Value = 23571.428
$ws.Range("N50").Value = -24821.428
$ws.Range("H51").Value = 20000
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").Value = $null
$ws.Range("H59").Value = 31666.334
$ws.Range("I59").Value = 34999
$ws.Range("K59").Value = 34999
$ws.Range("M59").Value = -33854
$ws.Range("H60").Value = 18920
$ws.Range("I60").Value = 18200
$ws.Range("J60").Value = 20000
$ws.Range("K60").Value = 18200
$ws.Range("L60").Value = 20000
$ws.Range("M60").Value = -17689
$ws.Range("N60").Value = -21022
$ws.Range("H61").Value = 20000
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").Value = $null
$ws.Range("H68").Value = 4999
$ws.Range("I68").Value = 4999
$ws.Range("K68").Value = 4999
$ws.Range("M68").Value = -4250
$ws.Range("H71").Value = 4999
$ws.Range("I71").Value = 4999
$ws.Range("K71").Value = 14997
$ws.Range("M71").Value = -11253

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 354.13333
$ws.Range("I12").Value = 8
$ws.Range("K12").Value = 24
$ws.Range("M12").Value = 149
$ws.Range("H18").Value = 1448.2727
$ws.Range("J18").Value = 1833.3334
$ws.Range("L18").Value = 5500.0002
$ws.Range("N18").Value = -5838.0002
$ws.Range("H34").Value = 682.8333
$ws.Range("I34").Value = 577.44446
$ws.Range("J34").Value = 999
$ws.Range("K34").Value = 1732.33338
$ws.Range("L34").Value = 2997
$ws.Range("M34").Value = -1648.33338
$ws.Range("N34").Value = -3165
$ws.Range("H39").Value = 1500
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").Value = $null
$ws.Range("H55").Value = 2833.3333
$ws.Range("I55").Value = 2000
$ws.Range("K55").Value = 6000
$ws.Range("M55").Value = -5823

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 12301.066
$ws.Range("I43").Value = 11319.637
$ws.Range("K43").Value = 11319.637
$ws.Range("M43").Value = -11168.637
$ws.Range("H46").Value = 8499.833000000001
$ws.Range("H57").Value = 12860.833
$ws.Range("H80").Value = 7500.143
$ws.Range("I80").Value = 3500
$ws.Range("J80").Value = 8166.8335
$ws.Range("K80").Value = 3500
$ws.Range("L80").Value = 8166.8335
$ws.Range("M80").Value = -2502
$ws.Range("N80").Value = -10162.8335
$ws.Range("H83").Value = 7500.143
$ws.Range("I83").Value = 3500
$ws.Range("J83").Value = 8166.8335
$ws.Range("K83").Value = 17500
$ws.Range("L83").Value = 40834.1675
$ws.Range("M83").Value = -12508
$ws.Range("N83").Value = -50818.1675

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 940.1579
$ws.Range("I82").Value = 1049.9231
$ws.Range("J82").Value = 702.3333
$ws.Range("K82").Value = 1049.9231
$ws.Range("L82").Value = 702.3333
$ws.Range("M82").Value = -688.9231
$ws.Range("N82").Value = -1424.3333
$ws.Range("H85").Value = 940.1579
$ws.Range("I85").Value = 1049.9231
$ws.Range("J85").Value = 702.3333
$ws.Range("K85").Value = 1049.9231
$ws.Range("L85").Value = 702.3333
$ws.Range("M85").Value = 198.0769
$ws.Range("N85").Value = -3198.3333
$ws.Range("H132").Value = 2281.889
$ws.Range("J132").Value = 2762.6667
$ws.Range("L132").Value = 8288.000100000001
$ws.Range("N132").Value = -13348.0001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 24939.541
$ws.Range("J54").Value = 31537.691
$ws.Range("L54").Value = 31537.691
$ws.Range("N54").Value = -32577.691
$ws.Range("H81").Value = 415.66666
$ws.Range("I81").Value = 415.66666
$ws.Range("K81").Value = 831.33332
$ws.Range("M81").Value = 229.66668
$ws.Range("H84").Value = 415.66666
$ws.Range("I84").Value = 415.66666
$ws.Range("K84").Value = 4156.6666
$ws.Range("M84").Value = 1147.3334
